$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I7 (Mlstr_harmo::comment) - updated to match Franzi's wording
$ws.Range("J7").Copy()
$ws.Range("I7").PasteSpecial(-4122)
$ws.Range("I7").Value = "non-smokers as inclusion criterion, but we cannot be sure, if they are real never smokers or former smokers"
$ws.Range("I7").Interior.ColorIndex = -4142

# K7 (Mlstr_harmo::status_detail) - changed to proximate
$ws.Range("K7").Value = "proximate"

# Row 7 no longer needs the tall custom height
$ws.Rows(7).AutoFit()

# Update active selection to I7
$ws.Range("I7").Select() | Out-Null
